# Updates crypto price (D) and 1h volume change (E) columns for rows 2-51.
# D-column values that parse as plain numbers must be round-tripped through a
# scratch cell formatted as Text and PasteSpecial(xlPasteValues) so Excel keeps
# storing them as text (matching the source data) instead of silently coercing
# them into numeric cells (which would drop significant trailing zeros, e.g.
# "1.00" -> 1, "17.40" -> 17.4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("ZZ1")

$ws.Range("D2").Value = '59.116.85'
$ws.Range("E2").Value = '  +0.64%  '

$ws.Range("D3").Value = '2.498.03'
$ws.Range("E3").Value = '  +0.85%  '

$scratch.NumberFormat = "@"
$scratch.Value = '1.00'
$scratch.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E4").Value = '  -0.12%  '

$scratch.NumberFormat = "@"
$scratch.Value = '536.64'
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E5").Value = '  +0.33%  '

$scratch.NumberFormat = "@"
$scratch.Value = '136.22'
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E6").Value = '  -0.46%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.999'
$scratch.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E7").Value = '  -0.31%  '

$ws.Range("E8").Value = '  +1.35%  '

$ws.Range("D9").Value = '2.518.27'
$ws.Range("E9").Value = '  +1.78%  '

$ws.Range("E10").Value = '  +1.66%  '

$ws.Range("E11").Value = '  -2.36%  '

$scratch.NumberFormat = "@"
$scratch.Value = '5.34'
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E12").Value = '  -1.17%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.348'
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E13").Value = '  +0.55%  '

$ws.Range("D14").Value = '2.947.02'
$ws.Range("E14").Value = '  +0.57%  '

$scratch.NumberFormat = "@"
$scratch.Value = '23.00'
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E15").Value = '  +1.24%  '

$ws.Range("D16").Value = '58.893.43'
$ws.Range("E16").Value = '  +0.34%  '

$ws.Range("E17").Value = '  +0.21%  '

$ws.Range("D18").Value = '2.516.48'
$ws.Range("E18").Value = '  +1.42%  '

$scratch.NumberFormat = "@"
$scratch.Value = '11.07'
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E19").Value = '  +2.19%  '

$scratch.NumberFormat = "@"
$scratch.Value = '4.26'
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E20").Value = '  +1.36%  '

$scratch.NumberFormat = "@"
$scratch.Value = '323.17'
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E21").Value = '  +0.58%  '

$ws.Range("E22").Value = '  +0.08%  '

$ws.Range("E23").Value = '  +3.74%  '

$scratch.NumberFormat = "@"
$scratch.Value = '65.07'
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E24").Value = '  +3.87%  '

$ws.Range("E25").Value = '  +2.44%  '

$ws.Range("E26").Value = '  -1.19%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.997'
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E27").Value = '  +0.19%  '

$ws.Range("E28").Value = '  -0.79%  '

$ws.Range("D29").Value = '0.0₃0768'
$ws.Range("E29").Value = '  +0.45%  '

$scratch.NumberFormat = "@"
$scratch.Value = '6.64'
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E30").Value = '  +0.88%  '

$ws.Range("E31").Value = '  -0.83%  '

$scratch.NumberFormat = "@"
$scratch.Value = '170.60'
$scratch.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E32").Value = '  +3.77%  '

$ws.Range("E33").Value = '  +9.27%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.998'
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E34").Value = '  -0.12%  '

$ws.Range("E35").Value = '  +1.44%  '

$ws.Range("E36").Value = '  +0.32%  '

$scratch.NumberFormat = "@"
$scratch.Value = '4.07'
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E37").Value = '  +0.42%  '

$scratch.NumberFormat = "@"
$scratch.Value = '1.54'
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E38").Value = '  -0.26%  '

$scratch.NumberFormat = "@"
$scratch.Value = '36.91'
$scratch.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E39").Value = '  +0.78%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.813'
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E40").Value = '  +3.28%  '

$ws.Range("E41").Value = '  +0.33%  '

$scratch.NumberFormat = "@"
$scratch.Value = '284.47'
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E42").Value = '  +3.50%  '

$scratch.NumberFormat = "@"
$scratch.Value = '5.07'
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E43").Value = '  -0.25%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.996'
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E44").Value = '  -0.71%  '

$ws.Range("E45").Value = '  +3.16%  '

$scratch.NumberFormat = "@"
$scratch.Value = '130.29'
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E46").Value = '  +5.99%  '

$scratch.NumberFormat = "@"
$scratch.Value = '10.87'
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E47").Value = '  -0.22%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.0923'
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E48").Value = '  -0.86%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.0504'
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E49").Value = '  -0.15%  '

$ws.Range("E50").Value = '  -0.15%  '

$scratch.NumberFormat = "@"
$scratch.Value = '17.40'
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E51").Value = '  +0.64%  '
